$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 with new values
$ws.Range("A2").Value = "0vJhfpLl"
$ws.Range("B2").Value = "27/11/2024"
$ws.Range("C2").Value = "09:45"
$ws.Range("D2").Value = "URUGUAY - PRIMERA DIVISION"
$ws.Range("E2").Value = "Fenix"
$ws.Range("F2").Value = "Maldonado"
$ws.Range("G2").Value = 1.73
$ws.Range("H2").Value = 3.5
$ws.Range("I2").Value = 5
$ws.Range("J2").Value = 2.4
$ws.Range("K2").Value = 2.05
$ws.Range("L2").Value = 5.5
$ws.Range("M2").Value = 1.07
$ws.Range("N2").Value = 9
$ws.Range("O2").Value = 1.36
$ws.Range("P2").Value = 3
$ws.Range("Q2").Value = 2.2
$ws.Range("R2").Value = 1.65
$ws.Range("S2").Value = 1.5
$ws.Range("T2").Value = 2.5
$ws.Range("U2").Value = 2.1
$ws.Range("V2").Value = 1.67
$ws.Range("W2").Value = 6
$ws.Range("X2").Value = 7.5
$ws.Range("Y2").Value = 9
$ws.Range("Z2").Value = 13
$ws.Range("AA2").Value = 17
$ws.Range("AB2").Value = 34
$ws.Range("AC2").Value = 8
$ws.Range("AD2").Value = 7
$ws.Range("AE2").Value = 19
$ws.Range("AF2").Value = 67
$ws.Range("AG2").Value = 1000
$ws.Range("AH2").Value = 11
$ws.Range("AI2").Value = 23
$ws.Range("AJ2").Value = 17
$ws.Range("AK2").Value = 51
$ws.Range("AL2").Value = 41
$ws.Range("AM2").Value = 51
$ws.Range("AN2").Value = 3.6
$ws.Range("AO2").Value = 9.5
$ws.Range("AP2").Value = 23
$ws.Range("AQ2").Value = 34
$ws.Range("AR2").Value = 51
$ws.Range("AS2").Value = 201
$ws.Range("AT2").Value = 2.5
$ws.Range("AU2").Value = 9
$ws.Range("AV2").Value = 67
$ws.Range("AW2").Value = 6.5
$ws.Range("AX2").Value = 29
$ws.Range("AY2").Value = 41
$ws.Range("AZ2").Value = 101
$ws.Range("BA2").Value = 151
$ws.Range("BB2").Value = 351

# Delete row 3 entirely (shift cells up)
$ws.Rows("3:3").Delete()
